$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 111756162
$ws.Range("Q8").Value = 453923
$ws.Range("R8").Value = 7073958

# Row 9
$ws.Range("A9").Value = 111756156
$ws.Range("B9").Value = 89423
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 5432
$ws.Range("F9").Value = 'Granticka'
$ws.Range("G9").Value = 'Porodaedalea chrysoloma'
$ws.Range("H9").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q9").Value = 453978
$ws.Range("R9").Value = 7073813

# Row 10
$ws.Range("A10").Value = 111756172
$ws.Range("B10").Value = 85715
$ws.Range("E10").Value = 510
$ws.Range("F10").Value = 'Doftskinn'
$ws.Range("G10").Value = 'Cystostereum murrayi'
$ws.Range("H10").Value = '(Berk. & M.A. Curtis.) Pouzar'
$ws.Range("Q10").Value = 453939
$ws.Range("R10").Value = 7073959

# Row 11
$ws.Range("A11").Value = 111756167
$ws.Range("B11").Value = 77515
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = 'Garnlav'
$ws.Range("G11").Value = 'Alectoria sarmentosa'
$ws.Range("H11").Value = '(Ach.) Ach.'
$ws.Range("Q11").Value = 454003
$ws.Range("R11").Value = 7073638

# Row 12
$ws.Range("A12").Value = 111756139
$ws.Range("B12").Value = 89405
$ws.Range("E12").Value = 1202
$ws.Range("F12").Value = 'Ullticka'
$ws.Range("G12").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H12").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q12").Value = 453693
$ws.Range("R12").Value = 7074032

# Row 13
$ws.Range("A13").Value = 111756141
$ws.Range("B13").Value = 89405
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = 'Ullticka'
$ws.Range("G13").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H13").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q13").Value = 453610
$ws.Range("R13").Value = 7074087

# Row 14
$ws.Range("A14").Value = 111756155
$ws.Range("B14").Value = 89423
$ws.Range("E14").Value = 5432
$ws.Range("F14").Value = 'Granticka'
$ws.Range("G14").Value = 'Porodaedalea chrysoloma'
$ws.Range("H14").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q14").Value = 453863
$ws.Range("R14").Value = 7073965

# Row 15
$ws.Range("A15").Value = 111756158
$ws.Range("B15").Value = 89423
$ws.Range("E15").Value = 5432
$ws.Range("F15").Value = 'Granticka'
$ws.Range("G15").Value = 'Porodaedalea chrysoloma'
$ws.Range("H15").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q15").Value = 454003
$ws.Range("R15").Value = 7073783

# Row 16
$ws.Range("A16").Value = 111756153
$ws.Range("B16").Value = 96674
$ws.Range("E16").Value = 219880
$ws.Range("F16").Value = 'Kransrams'
$ws.Range("G16").Value = 'Polygonatum verticillatum'
$ws.Range("H16").Value = '(L.) All.'
$ws.Range("Q16").Value = 453708
$ws.Range("R16").Value = 7073722

# Row 18
$ws.Range("A18").Value = 111756170
$ws.Range("B18").Value = 96265
$ws.Range("D18").Value = 'LC'
$ws.Range("E18").Value = 219790
$ws.Range("F18").Value = 'Fläcknycklar'
$ws.Range("G18").Value = 'Dactylorhiza maculata'
$ws.Range("H18").Value = '(L.) Soó'
$ws.Range("Q18").Value = 453739
$ws.Range("R18").Value = 7073724

# Row 19
$ws.Range("A19").Value = 111756150
$ws.Range("B19").Value = 95532
$ws.Range("E19").Value = 221945
$ws.Range("F19").Value = 'Revlummer'
$ws.Range("G19").Value = 'Lycopodium annotinum'
$ws.Range("H19").Value = 'L.'
$ws.Range("Q19").Value = 453976
$ws.Range("R19").Value = 7073812

# Row 20
$ws.Range("A20").Value = 111756140
$ws.Range("B20").Value = 89405
$ws.Range("E20").Value = 1202
$ws.Range("F20").Value = 'Ullticka'
$ws.Range("G20").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H20").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q20").Value = 453821
$ws.Range("R20").Value = 7074037

# Row 21
$ws.Range("A21").Value = 111756161
$ws.Range("B21").Value = 77515
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = 'Garnlav'
$ws.Range("G21").Value = 'Alectoria sarmentosa'
$ws.Range("H21").Value = '(Ach.) Ach.'
$ws.Range("Q21").Value = 453723
$ws.Range("R21").Value = 7074070

# Row 22
$ws.Range("A22").Value = 111756148
$ws.Range("B22").Value = 96266
$ws.Range("D22").Value = 'LC'
$ws.Range("E22").Value = 223591
$ws.Range("F22").Value = 'Skogsnycklar'
$ws.Range("G22").Value = 'Dactylorhiza maculata subsp. fuchsii'
$ws.Range("H22").Value = '(Druce) Hyl.'
$ws.Range("Q22").Value = 453747
$ws.Range("R22").Value = 7073851

# Row 23
$ws.Range("A23").Value = 111756151
$ws.Range("B23").Value = 95532
$ws.Range("D23").Value = 'LC'
$ws.Range("E23").Value = 221945
$ws.Range("F23").Value = 'Revlummer'
$ws.Range("G23").Value = 'Lycopodium annotinum'
$ws.Range("H23").Value = 'L.'
$ws.Range("Q23").Value = 453609
$ws.Range("R23").Value = 7074131

# Row 24
$ws.Range("A24").Value = 111756142
$ws.Range("B24").Value = 90087
$ws.Range("E24").Value = 3298
$ws.Range("F24").Value = 'Trådticka'
$ws.Range("G24").Value = 'Climacocystis borealis'
$ws.Range("H24").Value = '(Fr.) Kotl. & Pouzar'
$ws.Range("Q24").Value = 454003
$ws.Range("R24").Value = 7073638

# Row 25
$ws.Range("A25").Value = 111756143
$ws.Range("B25").Value = 90087
$ws.Range("D25").Value = 'LC'
$ws.Range("E25").Value = 3298
$ws.Range("F25").Value = 'Trådticka'
$ws.Range("G25").Value = 'Climacocystis borealis'
$ws.Range("H25").Value = '(Fr.) Kotl. & Pouzar'
$ws.Range("Q25").Value = 453951
$ws.Range("R25").Value = 7073592

# Row 26
$ws.Range("A26").Value = 111756157
$ws.Range("B26").Value = 89423
$ws.Range("D26").Value = 'NT'
$ws.Range("E26").Value = 5432
$ws.Range("F26").Value = 'Granticka'
$ws.Range("G26").Value = 'Porodaedalea chrysoloma'
$ws.Range("H26").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q26").Value = 453982
$ws.Range("R26").Value = 7073807

# Row 28
$ws.Range("A28").Value = 111756169
$ws.Range("B28").Value = 77515
$ws.Range("D28").Value = 'NT'
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = 'Garnlav'
$ws.Range("G28").Value = 'Alectoria sarmentosa'
$ws.Range("H28").Value = '(Ach.) Ach.'
$ws.Range("Q28").Value = 453910
$ws.Range("R28").Value = 7073654

# Row 29
$ws.Range("A29").Value = 111756159
$ws.Range("B29").Value = 89423
$ws.Range("E29").Value = 5432
$ws.Range("F29").Value = 'Granticka'
$ws.Range("G29").Value = 'Porodaedalea chrysoloma'
$ws.Range("H29").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q29").Value = 453621
$ws.Range("R29").Value = 7073984

# Row 30
$ws.Range("A30").Value = 111756147
$ws.Range("B30").Value = 89425
$ws.Range("E30").Value = 5442
$ws.Range("F30").Value = 'Tallticka'
$ws.Range("G30").Value = 'Porodaedalea pini'
$ws.Range("H30").Value = '(Brot.) Murrill'
$ws.Range("Q30").Value = 453989
$ws.Range("R30").Value = 7073710

# Row 31
$ws.Range("A31").Value = 111756168
$ws.Range("B31").Value = 77515
$ws.Range("E31").Value = 6425
$ws.Range("F31").Value = 'Garnlav'
$ws.Range("G31").Value = 'Alectoria sarmentosa'
$ws.Range("H31").Value = '(Ach.) Ach.'
$ws.Range("Q31").Value = 453959
$ws.Range("R31").Value = 7073596

# Row 32
$ws.Range("A32").Value = 111756160
$ws.Range("B32").Value = 77515
$ws.Range("D32").Value = 'NT'
$ws.Range("E32").Value = 6425
$ws.Range("F32").Value = 'Garnlav'
$ws.Range("G32").Value = 'Alectoria sarmentosa'
$ws.Range("H32").Value = '(Ach.) Ach.'
$ws.Range("Q32").Value = 453816
$ws.Range("R32").Value = 7073870

# Row 33
$ws.Range("A33").Value = 111756165
$ws.Range("B33").Value = 77515
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = 'Garnlav'
$ws.Range("G33").Value = 'Alectoria sarmentosa'
$ws.Range("H33").Value = '(Ach.) Ach.'
$ws.Range("Q33").Value = 453984
$ws.Range("R33").Value = 7073751

# Row 34
$ws.Range("A34").Value = 111756164
$ws.Range("B34").Value = 77515
$ws.Range("E34").Value = 6425
$ws.Range("F34").Value = 'Garnlav'
$ws.Range("G34").Value = 'Alectoria sarmentosa'
$ws.Range("H34").Value = '(Ach.) Ach.'
$ws.Range("Q34").Value = 453971
$ws.Range("R34").Value = 7073820

# Row 35
$ws.Range("A35").Value = 111756163
$ws.Range("B35").Value = 77515
$ws.Range("E35").Value = 6425
$ws.Range("F35").Value = 'Garnlav'
$ws.Range("G35").Value = 'Alectoria sarmentosa'
$ws.Range("H35").Value = '(Ach.) Ach.'
$ws.Range("Q35").Value = 453956
$ws.Range("R35").Value = 7073946
